$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update time_taken (column F) timestamps on the "data" sheet ---
$dataSheet.Range("F2").Value = "2021-10-05 14:19:26.667988"
$dataSheet.Range("F3").Value = "2021-10-05 14:19:26.667996"
$dataSheet.Range("F4").Value = "2021-10-05 14:19:26.668000"
$dataSheet.Range("F5").Value = "2021-10-05 14:19:26.668002"
$dataSheet.Range("F6").Value = "2021-10-05 14:19:26.668005"
$dataSheet.Range("F7").Value = "2021-10-05 14:19:26.668008"
$dataSheet.Range("F8").Value = "2021-10-05 14:19:26.668010"
$dataSheet.Range("F9").Value = "2021-10-05 14:19:26.668012"
$dataSheet.Range("F10").Value = "2021-10-05 14:19:26.668015"
$dataSheet.Range("F11").Value = "2021-10-05 14:19:26.668018"
$dataSheet.Range("F12").Value = "2021-10-05 14:19:26.668021"
$dataSheet.Range("F13").Value = "2021-10-05 14:19:26.668023"
$dataSheet.Range("F14").Value = "2021-10-05 14:19:26.668026"
$dataSheet.Range("F15").Value = "2021-10-05 14:19:26.668028"
$dataSheet.Range("F16").Value = "2021-10-05 14:19:26.668030"
$dataSheet.Range("F17").Value = "2021-10-05 14:19:26.668033"
$dataSheet.Range("F18").Value = "2021-10-05 14:19:26.668036"
$dataSheet.Range("F19").Value = "2021-10-05 14:19:26.668039"
$dataSheet.Range("F20").Value = "2021-10-05 14:19:26.668041"
$dataSheet.Range("F21").Value = "2021-10-05 14:19:26.668044"
$dataSheet.Range("F22").Value = "2021-10-05 14:19:26.668046"
$dataSheet.Range("F23").Value = "2021-10-05 14:19:26.668048"
$dataSheet.Range("F24").Value = "2021-10-05 14:19:26.668051"
$dataSheet.Range("F25").Value = "2021-10-05 14:19:26.668053"
$dataSheet.Range("F26").Value = "2021-10-05 14:19:26.668056"
$dataSheet.Range("F27").Value = "2021-10-05 14:19:26.668059"
$dataSheet.Range("F28").Value = "2021-10-05 14:19:26.668061"
$dataSheet.Range("F29").Value = "2021-10-05 14:19:26.668064"
$dataSheet.Range("F30").Value = "2021-10-05 14:19:26.668066"
$dataSheet.Range("F31").Value = "2021-10-05 14:19:26.668069"
$dataSheet.Range("F32").Value = "2021-10-05 14:19:26.668071"
$dataSheet.Range("F33").Value = "2021-10-05 14:19:26.668073"
$dataSheet.Range("F34").Value = "2021-10-05 14:19:26.668076"
$dataSheet.Range("F35").Value = "2021-10-05 14:19:26.668079"
$dataSheet.Range("F36").Value = "2021-10-05 14:19:26.668081"
$dataSheet.Range("F37").Value = "2021-10-05 14:19:26.668084"
$dataSheet.Range("F38").Value = "2021-10-05 14:19:26.668086"
$dataSheet.Range("F39").Value = "2021-10-05 14:19:26.668089"
$dataSheet.Range("F40").Value = "2021-10-05 14:19:26.668091"
$dataSheet.Range("F41").Value = "2021-10-05 14:19:26.668094"
$dataSheet.Range("F42").Value = "2021-10-05 14:19:26.668096"
$dataSheet.Range("F43").Value = "2021-10-05 14:19:26.668099"
$dataSheet.Range("F44").Value = "2021-10-05 14:19:26.668101"
$dataSheet.Range("F45").Value = "2021-10-05 14:19:26.668104"
$dataSheet.Range("F46").Value = "2021-10-05 14:19:26.668106"
$dataSheet.Range("F47").Value = "2021-10-05 14:19:26.668109"
$dataSheet.Range("F48").Value = "2021-10-05 14:19:26.668111"
$dataSheet.Range("F49").Value = "2021-10-05 14:19:26.668114"
$dataSheet.Range("F50").Value = "2021-10-05 14:19:26.668116"
$dataSheet.Range("F51").Value = "2021-10-05 14:19:26.668119"
$dataSheet.Range("F52").Value = "2021-10-05 14:19:26.668121"
$dataSheet.Range("F53").Value = "2021-10-05 14:19:26.668124"
$dataSheet.Range("F54").Value = "2021-10-05 14:19:26.668127"
$dataSheet.Range("F55").Value = "2021-10-05 14:19:26.668129"
$dataSheet.Range("F56").Value = "2021-10-05 14:19:26.668132"
$dataSheet.Range("F57").Value = "2021-10-05 14:19:26.668134"
$dataSheet.Range("F58").Value = "2021-10-05 14:19:26.668137"
$dataSheet.Range("F59").Value = "2021-10-05 14:19:26.668139"
$dataSheet.Range("F60").Value = "2021-10-05 14:19:26.668142"
$dataSheet.Range("F61").Value = "2021-10-05 14:19:26.668144"
$dataSheet.Range("F62").Value = "2021-10-05 14:19:26.668147"
$dataSheet.Range("F63").Value = "2021-10-05 14:19:26.668149"
$dataSheet.Range("F64").Value = "2021-10-05 14:19:26.668152"
$dataSheet.Range("F65").Value = "2021-10-05 14:19:26.668154"
$dataSheet.Range("F66").Value = "2021-10-05 14:19:26.668158"
$dataSheet.Range("F67").Value = "2021-10-05 14:19:26.668160"
$dataSheet.Range("F68").Value = "2021-10-05 14:19:26.668163"
$dataSheet.Range("F69").Value = "2021-10-05 14:19:26.668165"
$dataSheet.Range("F70").Value = "2021-10-05 14:19:26.668168"
$dataSheet.Range("F71").Value = "2021-10-05 14:19:26.668170"
$dataSheet.Range("F72").Value = "2021-10-05 14:19:26.668173"
$dataSheet.Range("F73").Value = "2021-10-05 14:19:26.668175"
$dataSheet.Range("F74").Value = "2021-10-05 14:19:26.668178"
$dataSheet.Range("F75").Value = "2021-10-05 14:19:26.668180"
$dataSheet.Range("F76").Value = "2021-10-05 14:19:26.668183"
$dataSheet.Range("F77").Value = "2021-10-05 14:19:26.668185"
$dataSheet.Range("F78").Value = "2021-10-05 14:19:26.668190"
$dataSheet.Range("F79").Value = "2021-10-05 14:19:26.668193"
$dataSheet.Range("F80").Value = "2021-10-05 14:19:26.668196"
$dataSheet.Range("F81").Value = "2021-10-05 14:19:26.668198"
$dataSheet.Range("F82").Value = "2021-10-05 14:19:26.668201"
$dataSheet.Range("F83").Value = "2021-10-05 14:19:26.668203"
$dataSheet.Range("F84").Value = "2021-10-05 14:19:26.668206"
$dataSheet.Range("F85").Value = "2021-10-05 14:19:26.668208"
$dataSheet.Range("F86").Value = "2021-10-05 14:19:26.668211"
$dataSheet.Range("F87").Value = "2021-10-05 14:19:26.668214"
$dataSheet.Range("F88").Value = "2021-10-05 14:19:26.668216"
$dataSheet.Range("F89").Value = "2021-10-05 14:19:26.668219"
$dataSheet.Range("F90").Value = "2021-10-05 14:19:26.668221"
$dataSheet.Range("F91").Value = "2021-10-05 14:19:26.668224"
$dataSheet.Range("F92").Value = "2021-10-05 14:19:26.668227"
$dataSheet.Range("F93").Value = "2021-10-05 14:19:26.668229"
$dataSheet.Range("F94").Value = "2021-10-05 14:19:26.668233"
$dataSheet.Range("F95").Value = "2021-10-05 14:19:26.668236"
$dataSheet.Range("F96").Value = "2021-10-05 14:19:26.668239"
$dataSheet.Range("F97").Value = "2021-10-05 14:19:26.668241"
$dataSheet.Range("F98").Value = "2021-10-05 14:19:26.668244"
$dataSheet.Range("F99").Value = "2021-10-05 14:19:26.668246"
$dataSheet.Range("F100").Value = "2021-10-05 14:19:26.668249"
$dataSheet.Range("F101").Value = "2021-10-05 14:19:26.668251"
$dataSheet.Range("F102").Value = "2021-10-05 14:19:26.668254"
$dataSheet.Range("F103").Value = "2021-10-05 14:19:26.668256"
$dataSheet.Range("F104").Value = "2021-10-05 14:19:26.668259"
$dataSheet.Range("F105").Value = "2021-10-05 14:19:26.668261"
$dataSheet.Range("F106").Value = "2021-10-05 14:19:26.668264"
$dataSheet.Range("F107").Value = "2021-10-05 14:19:26.668266"
$dataSheet.Range("F108").Value = "2021-10-05 14:19:26.668269"
$dataSheet.Range("F109").Value = "2021-10-05 14:19:26.668272"
$dataSheet.Range("F110").Value = "2021-10-05 14:19:26.668276"
$dataSheet.Range("F111").Value = "2021-10-05 14:19:26.668279"
$dataSheet.Range("F112").Value = "2021-10-05 14:19:26.668282"
$dataSheet.Range("F113").Value = "2021-10-05 14:19:26.668284"
$dataSheet.Range("F114").Value = "2021-10-05 14:19:26.668287"
$dataSheet.Range("F115").Value = "2021-10-05 14:19:26.668289"
$dataSheet.Range("F116").Value = "2021-10-05 14:19:26.668292"
$dataSheet.Range("F117").Value = "2021-10-05 14:19:26.668294"
$dataSheet.Range("F118").Value = "2021-10-05 14:19:26.668297"
$dataSheet.Range("F119").Value = "2021-10-05 14:19:26.668299"
$dataSheet.Range("F120").Value = "2021-10-05 14:19:26.668302"
$dataSheet.Range("F121").Value = "2021-10-05 14:19:26.668304"
$dataSheet.Range("F122").Value = "2021-10-05 14:19:26.668307"
$dataSheet.Range("F123").Value = "2021-10-05 14:19:26.668309"
$dataSheet.Range("F124").Value = "2021-10-05 14:19:26.668312"
$dataSheet.Range("F125").Value = "2021-10-05 14:19:26.668315"
$dataSheet.Range("F126").Value = "2021-10-05 14:19:26.668317"
$dataSheet.Range("F127").Value = "2021-10-05 14:19:26.668320"
$dataSheet.Range("F128").Value = "2021-10-05 14:19:26.668322"
$dataSheet.Range("F129").Value = "2021-10-05 14:19:26.668324"
$dataSheet.Range("F130").Value = "2021-10-05 14:19:26.668329"
$dataSheet.Range("F131").Value = "2021-10-05 14:19:26.668332"
$dataSheet.Range("F132").Value = "2021-10-05 14:19:26.668334"
$dataSheet.Range("F133").Value = "2021-10-05 14:19:26.668337"
$dataSheet.Range("F134").Value = "2021-10-05 14:19:26.668339"
$dataSheet.Range("F135").Value = "2021-10-05 14:19:26.668342"
$dataSheet.Range("F136").Value = "2021-10-05 14:19:26.668344"
$dataSheet.Range("F137").Value = "2021-10-05 14:19:26.668347"
$dataSheet.Range("F138").Value = "2021-10-05 14:19:26.668349"
$dataSheet.Range("F139").Value = "2021-10-05 14:19:26.668352"
$dataSheet.Range("F140").Value = "2021-10-05 14:19:26.668354"
$dataSheet.Range("F141").Value = "2021-10-05 14:19:26.668357"
$dataSheet.Range("F142").Value = "2021-10-05 14:19:26.668360"
$dataSheet.Range("F143").Value = "2021-10-05 14:19:26.668362"
$dataSheet.Range("F144").Value = "2021-10-05 14:19:26.668365"
$dataSheet.Range("F145").Value = "2021-10-05 14:19:26.668367"
$dataSheet.Range("F146").Value = "2021-10-05 14:19:26.668370"
$dataSheet.Range("F147").Value = "2021-10-05 14:19:26.668373"
$dataSheet.Range("F148").Value = "2021-10-05 14:19:26.668375"
$dataSheet.Range("F149").Value = "2021-10-05 14:19:26.668377"
$dataSheet.Range("F150").Value = "2021-10-05 14:19:26.668380"
$dataSheet.Range("F151").Value = "2021-10-05 14:19:26.668383"
$dataSheet.Range("F152").Value = "2021-10-05 14:19:26.668385"
$dataSheet.Range("F153").Value = "2021-10-05 14:19:26.668388"
$dataSheet.Range("F154").Value = "2021-10-05 14:19:26.668391"
$dataSheet.Range("F155").Value = "2021-10-05 14:19:26.668393"
$dataSheet.Range("F156").Value = "2021-10-05 14:19:26.668396"
$dataSheet.Range("F157").Value = "2021-10-05 14:19:26.668398"
$dataSheet.Range("F158").Value = "2021-10-05 14:19:26.668401"
$dataSheet.Range("F159").Value = "2021-10-05 14:19:26.668404"
$dataSheet.Range("F160").Value = "2021-10-05 14:19:26.668406"
$dataSheet.Range("F161").Value = "2021-10-05 14:19:26.668409"
$dataSheet.Range("F162").Value = "2021-10-05 14:19:26.668412"
$dataSheet.Range("F163").Value = "2021-10-05 14:19:26.668414"
$dataSheet.Range("F164").Value = "2021-10-05 14:19:26.668417"
$dataSheet.Range("F165").Value = "2021-10-05 14:19:26.668420"
$dataSheet.Range("F166").Value = "2021-10-05 14:19:26.668422"
$dataSheet.Range("F167").Value = "2021-10-05 14:19:26.668424"
$dataSheet.Range("F168").Value = "2021-10-05 14:19:26.668427"
$dataSheet.Range("F169").Value = "2021-10-05 14:19:26.668430"
$dataSheet.Range("F170").Value = "2021-10-05 14:19:26.668432"
$dataSheet.Range("F171").Value = "2021-10-05 14:19:26.668435"
$dataSheet.Range("F172").Value = "2021-10-05 14:19:26.668438"
$dataSheet.Range("F173").Value = "2021-10-05 14:19:26.668440"
$dataSheet.Range("F174").Value = "2021-10-05 14:19:26.668444"
$dataSheet.Range("F175").Value = "2021-10-05 14:19:26.668447"
$dataSheet.Range("F176").Value = "2021-10-05 14:19:26.668450"
$dataSheet.Range("F177").Value = "2021-10-05 14:19:26.668452"
$dataSheet.Range("F178").Value = "2021-10-05 14:19:26.668455"
$dataSheet.Range("F179").Value = "2021-10-05 14:19:26.668458"
$dataSheet.Range("F180").Value = "2021-10-05 14:19:26.668460"
$dataSheet.Range("F181").Value = "2021-10-05 14:19:26.668463"
$dataSheet.Range("F182").Value = "2021-10-05 14:19:26.668465"
$dataSheet.Range("F183").Value = "2021-10-05 14:19:26.668468"
$dataSheet.Range("F184").Value = "2021-10-05 14:19:26.668471"
$dataSheet.Range("F185").Value = "2021-10-05 14:19:26.668473"
$dataSheet.Range("F186").Value = "2021-10-05 14:19:26.668476"
$dataSheet.Range("F187").Value = "2021-10-05 14:19:26.668479"
$dataSheet.Range("F188").Value = "2021-10-05 14:19:26.668481"
$dataSheet.Range("F189").Value = "2021-10-05 14:19:26.668484"
$dataSheet.Range("F190").Value = "2021-10-05 14:19:26.668486"
$dataSheet.Range("F191").Value = "2021-10-05 14:19:26.668489"
$dataSheet.Range("F192").Value = "2021-10-05 14:19:26.668492"
$dataSheet.Range("F193").Value = "2021-10-05 14:19:26.668495"
$dataSheet.Range("F194").Value = "2021-10-05 14:19:26.668498"
$dataSheet.Range("F195").Value = "2021-10-05 14:19:26.668500"
$dataSheet.Range("F196").Value = "2021-10-05 14:19:26.668503"
$dataSheet.Range("F197").Value = "2021-10-05 14:19:26.668506"
$dataSheet.Range("F198").Value = "2021-10-05 14:19:26.668508"

# --- Add the new "metadata" worksheet right after "data" ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (bold, centered/top aligned, thin border - matches "data" sheet header style)
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$headerRange = $metaSheet.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("A2").Font.Bold = $true
$metaSheet.Range("A2").HorizontalAlignment = -4108
$metaSheet.Range("A2").VerticalAlignment = -4160
$metaSheet.Range("A2").Borders.LineStyle = 1

$metaSheet.Range("B2").Value = "Cataracts"
$metaSheet.Range("C2").Value = 230
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "2.84"
$metaSheet.Range("E2").NumberFormat = "@"
$metaSheet.Range("E2").Value = "2021-09-14T09:26:19.150844Z"
$metaSheet.Range("F2").NumberFormat = "@"
$metaSheet.Range("F2").Value = "2021-10-05 14:19:26.664681"
$metaSheet.Range("G2").NumberFormat = "@"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/230/?format=json"

# Keep "data" as the active sheet/tab (matches target activeTab="0")
$dataSheet.Activate()
